# Add a new "Calculadora" worksheet ahead of the existing "Dominio y Hosting"
# sheet, with a small table of hourly-rate / hours fields used for budgeting.

$wb = $excel.ActiveWorkbook

# Worksheets.Add() inserts the new sheet before the currently active sheet
# and makes it the active sheet - matching the target layout where
# "Calculadora" becomes sheet index 0 (sheetId 2, rId1) and
# "Dominio y Hosting" shifts to rId2.
$calc = $wb.Worksheets.Add()
$calc.Name = "Calculadora"

$calc.Range("A1").Value = "Precio hora"
$calc.Range("B1").Value = 400
$calc.Range("A2").Value = "Horas de diseño"
$calc.Range("A3").Value = "Horas de prototipado"
$calc.Range("A4").Value = "Horas de desarrollo"
$calc.Range("A5").Value = "Horas de prueba"
$calc.Range("A6").Value = "Horas de subida"

# Leave the cursor parked a couple of rows below the table, mirroring the
# saved selection state (A8) captured in the authored workbook.
[void]$calc.Range("A8").Select()
